# Add a new data row (row 2) to the foot-print database sheet, mirroring
# the "alona hoz" / "alona" / "alon1105" record added by the author.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: cells are written in this particular order (Username, then
# Full name, then Password) so that new shared-string entries land in the
# same order as the target workbook: "alona" (Username) is inserted into
# the shared string table first, then "alona hoz" (Full name), then
# "alon1105" (Password).
$ws.Range("C2").Value = "alona"
$ws.Range("A2").Value = "alona hoz"
$ws.Range("D2").Value = "alon1105"

$ws.Range("B2").Value = 329458418

$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 150
$ws.Range("G2").Value = 160
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 9
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0

# Column B ("ID_NUMBER") got auto-fitted to its content in the author's
# edit (9-digit id numbers); apply a matching custom width.
$ws.Columns.Item(2).ColumnWidth = 9

# The author's final selection (visible in the saved file) was cell O4.
$ws.Range("O4").Select()
